# Applies the "Updated symbol list" crypto-price refresh described by the diff.
# Source data cells are inline strings (coin names, URLs, price text, percent text),
# so every numeric-looking value is written as literal TEXT (NumberFormat "@")
# to avoid Excel auto-converting "0.1900" -> 0.19 or "-0.32%" -> a percentage number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "305.42"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.31%"
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.58"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.32%"
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.043"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.53%"
# Row 5
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.54%"
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.875"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-2.25%"
# Row 7
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.790"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.77%"
# Row 8
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9228"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.77%"
# Row 9
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1284"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-7.64%"
# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1900"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.10%"
# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09129"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.94%"
# Row 12
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03416"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-5.73%"
# Row 13
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09862"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.54%"
# Row 14
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001400"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.43%"
# Row 15
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006197"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "4.98%"
# Row 16
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.853"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "8.45%"
# Row 17
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.132"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.22%"
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.352"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "11.18%"
# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.37%"
# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "3.37%"
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.996"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.33%"
# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-7.95%"
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04422"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.39%"
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001234"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.96%"
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004887"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.12%"
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001251"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-19.89%"
# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "42.08%"
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01936"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-1.24%"
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05176"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "5.49%"
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007610"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.28%"
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01014"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "9.63%"
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1349"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.73%"
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002152"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2.39%"
# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-15.10%"
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006176"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-3.22%"
# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.01%"
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "64.97"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-0.38%"
# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.01%"
# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.01%"
